$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename the header cells B1:N1 so they begin with a non-numeric
# character: prefix the existing "1.a".."2.h" labels with "Q".
$ws.Range("B1").Value = "Q1.a"
$ws.Range("C1").Value = "Q1.b"
$ws.Range("D1").Value = "Q1.c"
$ws.Range("E1").Value = "Q1.d"
$ws.Range("F1").Value = "Q1.e"
$ws.Range("G1").Value = "Q2.a"
$ws.Range("H1").Value = "Q2.b"
$ws.Range("I1").Value = "Q2.c"
$ws.Range("J1").Value = "Q2.d"
$ws.Range("K1").Value = "Q2.e"
$ws.Range("L1").Value = "Q2.f"
$ws.Range("M1").Value = "Q2.g"
$ws.Range("N1").Value = "Q2.h"

# Leave the rest of the header row (A1, O1:R1) untouched.

# Move the view / selection the way the author left it: scrolled down so
# row 47 is the first visible row under the frozen header, with K57 as
# the active cell.
$ws.Activate()
$ws.Range("K57").Select()
$excel.ActiveWindow.ScrollRow = 47
